# Add 2022-Q3 data
# 1) Insert a new "2022-Q3" sheet right after the summary ("总计") sheet by
#    duplicating the existing "2022-Q2" sheet (so formatting/styles match
#    exactly) and overwriting its values with the 2022-Q3 figures.
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q3 and
#    shift the existing quarters down by one.

$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# --- Step 1: create the new "2022-Q3" sheet ------------------------------
$q2Sheet.Copy($null, $summarySheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

$fundData = @(
  @(0, "161128", "易方达标普信息科技指数（QDII-LOF）人民币",   "4.99", "91.96", "1.91", "0.0953", 7),
  @(1, "012868", "易方达标普信息科技指数（QDII-LOF）人民币 C", "4.99", "91.96", "1.91", "0.0953", 7),
  @(2, "003721", "易方达标普信息科技指数（QDII-LOF）美元A",    "4.84", "91.96", "1.91", "0.0924", 7),
  @(3, "012869", "易方达标普信息科技指数（QDII-LOF）美元 C",   "0.15", "91.96", "1.91", "0.0029", 7)
)

$r = 2
foreach ($row in $fundData) {
  $q3Sheet.Cells.Item($r, 1).Value = $row[0]
  $q3Sheet.Cells.Item($r, 2).Value = "'" + $row[1]
  $q3Sheet.Cells.Item($r, 3).Value = $row[2]
  $q3Sheet.Cells.Item($r, 4).Value = "'" + $row[3]
  $q3Sheet.Cells.Item($r, 5).Value = "'" + $row[4]
  $q3Sheet.Cells.Item($r, 6).Value = "'" + $row[5]
  $q3Sheet.Cells.Item($r, 7).Value = "'" + $row[6]
  $q3Sheet.Cells.Item($r, 8).Value = $row[7]
  $r = $r + 1
}

# --- Step 2: update the "总计" summary sheet ------------------------------
$summarySheet.Range("A2:D2").Insert()
$summarySheet.Range("B2:D2").ClearFormats()

$summarySheet.Cells.Item(2, 1).Value = 0
$summarySheet.Cells.Item(2, 2).Value = "2022-Q3"
$summarySheet.Cells.Item(2, 3).Value = 4
$summarySheet.Cells.Item(2, 4).Value = 0.29

# Match the index-column style (bold/centered/bordered) used by the other
# rows in column A.
$summarySheet.Range("A3").Copy()
$summarySheet.Range("A2").PasteSpecial(-4122)

# Re-number the index column (A) for the rows that shifted down.
for ($i = 0; $i -lt 5; $i++) {
  $row = 3 + $i
  $summarySheet.Cells.Item($row, 1).Value = $i + 1
}

$excel.CutCopyMode = $false

# Keep the originally-active sheet/tab selected, matching the source file.
$summarySheet.Activate()
